# texts.xlsx update: TouchGFX typography / translation tables
#
# - Typography sheet: fill in "Wildcard Characters" (G) and "Wildcard Ranges" (I)
#   for the "Default" typography row so it can render the new "-" and "0-9"
#   glyphs used by the new screens below.
# - Translation sheet: re-purpose / add rows for the new "boat selection" and
#   "it works" screens (SingleUseId2 .. SingleUseId7).

$wb = $excel.ActiveWorkbook

$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("G4").Value = "-"
$wsTypo.Range("I4").Value = "0-9"

$wsTr = $wb.Worksheets.Item("Translation")

# Row 4: was SingleUseId1 / Default / Center / LTR / "Connexion >"
#   -> becomes SingleUseId2 / Default / Center / LTR / "Sélection / du voilier"
$wsTr.Range("B4").Value = "SingleUseId2"
$wsTr.Range("C4").Value = "Default"
$wsTr.Range("D4").Value = "Center"
$wsTr.Range("E4").Value = "LTR"
$wsTr.Range("F4").Value = "Sélection`ndu voilier"

# Row 5: was SingleUseId2 / Default / Left / LTR / "Selection du voilier"
#   -> becomes SingleUseId3 / Canal / Center / LTR / "<value>"
$wsTr.Range("B5").Value = "SingleUseId3"
$wsTr.Range("C5").Value = "Canal"
$wsTr.Range("D5").Value = "Center"
$wsTr.Range("E5").Value = "LTR"
$wsTr.Range("F5").Value = "<value>"

# Row 6: was SingleUseId3 / Canal / Center / LTR / "<value>"
#   -> becomes SingleUseId4 / Canal / Left / LTR / "0"
$wsTr.Range("B6").Value = "SingleUseId4"
$wsTr.Range("C6").Value = "Canal"
$wsTr.Range("D6").Value = "Left"
$wsTr.Range("E6").Value = "LTR"
# force text storage ("0" must stay a label, not turn into the number 0)
$wsTr.Range("F6").NumberFormat = "@"
$wsTr.Range("F6").Value = "0"
$wsTr.Range("F6").Style = "Normal"

# Row 7: was SingleUseId4 / Canal / Left / LTR / "0"
#   -> becomes SingleUseId5 / Default / Center / LTR / "<value>%"
$wsTr.Range("B7").Value = "SingleUseId5"
$wsTr.Range("C7").Value = "Default"
$wsTr.Range("D7").Value = "Center"
$wsTr.Range("E7").Value = "LTR"
$wsTr.Range("F7").Value = "<value>%"

# Row 8 (new): SingleUseId6 / Default / Left / LTR / "0"
$wsTr.Range("B8").Value = "SingleUseId6"
$wsTr.Range("C8").Value = "Default"
$wsTr.Range("D8").Value = "Left"
$wsTr.Range("E8").Value = "LTR"
$wsTr.Range("F8").NumberFormat = "@"
$wsTr.Range("F8").Value = "0"
$wsTr.Range("F8").Style = "Normal"

# Row 9 (new): SingleUseId7 / Small / Left / LTR / "<value>"
$wsTr.Range("B9").Value = "SingleUseId7"
$wsTr.Range("C9").Value = "Small"
$wsTr.Range("D9").Value = "Left"
$wsTr.Range("E9").Value = "LTR"
$wsTr.Range("F9").Value = "<value>"
